$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy row formatting (alternating date banding) from template rows ---
# Rows 908-914 (2025-03-03, Monday) continue the banding after row 907 (style family 29/30/31),
# so they take the opposite family 9/10/11, matching template row 893.
[void]$ws.Range("A893:T893").Copy()
$ws.Range("A908:T914").PasteSpecial(-4122)

# Rows 915-923 (2025-03-04, Tuesday) take the opposite family again (29/30/31),
# matching template row 899.
[void]$ws.Range("A899:T899").Copy()
$ws.Range("A915:T923").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Step 2: explicit row height (matches ht="20" customHeight="1" on each new row) ---
$ws.Range("908:923").RowHeight = 20

# --- Step 3: fill in the new game rows (2025-03-03 and 2025-03-04 slates) ---
# Row 908: CHA @ GSW (Monday 45719)
$ws.Cells.Item(908,1).Value = 45719
$ws.Cells.Item(908,2).Value = "Monday"
$ws.Cells.Item(908,3).Value = "GSW"
$ws.Cells.Item(908,4).Value = "CHA"
$ws.Cells.Item(908,5).Value = 0
$ws.Cells.Item(908,6).Value = "Mark Lindsay"
$ws.Cells.Item(908,7).Value = "Andy Nagy"
$ws.Cells.Item(908,8).Value = "Michael Smith"
$ws.Cells.Item(908,9).Value = 221.5
$ws.Cells.Item(908,10).Value = 12.5
$ws.Cells.Item(908,11).Value = 119
$ws.Cells.Item(908,12).Value = 101
$ws.Cells.Item(908,13).Formula = "=K908+L908"
$ws.Cells.Item(908,14).Formula = "=(L908-K908)*-1"
$ws.Cells.Item(908,15).Value = 1
$ws.Cells.Item(908,16).Formula = "=IF(M908>I908,1,0)"
$ws.Cells.Item(908,17).Formula = "=IF(P908=1,(M908-I908), """")"
$ws.Cells.Item(908,18).Formula = "=IF(M908<I908, 1, 0)"
$ws.Cells.Item(908,19).Formula = "=IF(R908=1,(I908-M908),"""")"
$ws.Cells.Item(908,20).Formula = "=IF(M908=I908,1,0)"

# Row 909: PHI @ POR (Monday 45719)
$ws.Cells.Item(909,1).Value = 45719
$ws.Cells.Item(909,2).Value = "Monday"
$ws.Cells.Item(909,3).Value = "POR"
$ws.Cells.Item(909,4).Value = "PHI"
$ws.Cells.Item(909,5).Value = 0
$ws.Cells.Item(909,6).Value = "Sean Wright"
$ws.Cells.Item(909,7).Value = "Marat Kogut"
$ws.Cells.Item(909,8).Value = "Danielle Scott"
$ws.Cells.Item(909,9).Value = 222
$ws.Cells.Item(909,10).Value = -3.5
$ws.Cells.Item(909,11).Value = 119
$ws.Cells.Item(909,12).Value = 102
$ws.Cells.Item(909,13).Formula = "=K909+L909"
$ws.Cells.Item(909,14).Formula = "=(L909-K909)*-1"
$ws.Cells.Item(909,15).Value = 1
$ws.Cells.Item(909,16).Formula = "=IF(M909>I909,1,0)"
$ws.Cells.Item(909,17).Formula = "=IF(P909=1,(M909-I909), """")"
$ws.Cells.Item(909,18).Formula = "=IF(M909<I909, 1, 0)"
$ws.Cells.Item(909,19).Formula = "=IF(R909=1,(I909-M909),"""")"
$ws.Cells.Item(909,20).Formula = "=IF(M909=I909,1,0)"

# Row 910: MIA @ WAS (Monday 45719)
$ws.Cells.Item(910,1).Value = 45719
$ws.Cells.Item(910,2).Value = "Monday"
$ws.Cells.Item(910,3).Value = "WAS"
$ws.Cells.Item(910,4).Value = "MIA"
$ws.Cells.Item(910,5).Value = 0
$ws.Cells.Item(910,6).Value = "Curtis Blair"
$ws.Cells.Item(910,7).Value = "Phenizee Ransom"
$ws.Cells.Item(910,8).Value = "Derrick Collins"
$ws.Cells.Item(910,9).Value = 225.5
$ws.Cells.Item(910,10).Value = -12.5
$ws.Cells.Item(910,11).Value = 90
$ws.Cells.Item(910,12).Value = 106
$ws.Cells.Item(910,13).Formula = "=K910+L910"
$ws.Cells.Item(910,14).Formula = "=(L910-K910)*-1"
$ws.Cells.Item(910,15).Value = 1
$ws.Cells.Item(910,16).Formula = "=IF(M910>I910,1,0)"
$ws.Cells.Item(910,17).Formula = "=IF(P910=1,(M910-I910), """")"
$ws.Cells.Item(910,18).Formula = "=IF(M910<I910, 1, 0)"
$ws.Cells.Item(910,19).Formula = "=IF(R910=1,(I910-M910),"""")"
$ws.Cells.Item(910,20).Formula = "=IF(M910=I910,1,0)"

# Row 911: MEM @ ATL (Monday 45719)
$ws.Cells.Item(911,1).Value = 45719
$ws.Cells.Item(911,2).Value = "Monday"
$ws.Cells.Item(911,3).Value = "ATL"
$ws.Cells.Item(911,4).Value = "MEM"
$ws.Cells.Item(911,5).Value = 0
$ws.Cells.Item(911,6).Value = "John Goble"
$ws.Cells.Item(911,7).Value = "Ray Acosta"
$ws.Cells.Item(911,8).Value = "Jonathan Sterling"
$ws.Cells.Item(911,9).Value = 250
$ws.Cells.Item(911,10).Value = -6
$ws.Cells.Item(911,11).Value = 132
$ws.Cells.Item(911,12).Value = 130
$ws.Cells.Item(911,13).Formula = "=K911+L911"
$ws.Cells.Item(911,14).Formula = "=(L911-K911)*-1"
$ws.Cells.Item(911,15).Value = 1
$ws.Cells.Item(911,16).Formula = "=IF(M911>I911,1,0)"
$ws.Cells.Item(911,17).Formula = "=IF(P911=1,(M911-I911), """")"
$ws.Cells.Item(911,18).Formula = "=IF(M911<I911, 1, 0)"
$ws.Cells.Item(911,19).Formula = "=IF(R911=1,(I911-M911),"""")"
$ws.Cells.Item(911,20).Formula = "=IF(M911=I911,1,0)"

# Row 912: OKC @ HOU (Monday 45719)
$ws.Cells.Item(912,1).Value = 45719
$ws.Cells.Item(912,2).Value = "Monday"
$ws.Cells.Item(912,3).Value = "HOU"
$ws.Cells.Item(912,4).Value = "OKC"
$ws.Cells.Item(912,5).Value = 0
$ws.Cells.Item(912,6).Value = "Ed Malloy"
$ws.Cells.Item(912,7).Value = "Justin Van Duyne"
$ws.Cells.Item(912,8).Value = "Matt Myers"
$ws.Cells.Item(912,9).Value = 220.5
$ws.Cells.Item(912,10).Value = -11
$ws.Cells.Item(912,11).Value = 128
$ws.Cells.Item(912,12).Value = 137
$ws.Cells.Item(912,13).Formula = "=K912+L912"
$ws.Cells.Item(912,14).Formula = "=(L912-K912)*-1"
$ws.Cells.Item(912,15).Value = 1
$ws.Cells.Item(912,16).Formula = "=IF(M912>I912,1,0)"
$ws.Cells.Item(912,17).Formula = "=IF(P912=1,(M912-I912), """")"
$ws.Cells.Item(912,18).Formula = "=IF(M912<I912, 1, 0)"
$ws.Cells.Item(912,19).Formula = "=IF(R912=1,(I912-M912),"""")"
$ws.Cells.Item(912,20).Formula = "=IF(M912=I912,1,0)"

# Row 913: DAL @ SAC (Monday 45719)
$ws.Cells.Item(913,1).Value = 45719
$ws.Cells.Item(913,2).Value = "Monday"
$ws.Cells.Item(913,3).Value = "SAC"
$ws.Cells.Item(913,4).Value = "DAL"
$ws.Cells.Item(913,5).Value = 0
$ws.Cells.Item(913,6).Value = "Tyler Ford"
$ws.Cells.Item(913,7).Value = "Jason Goldenberg"
$ws.Cells.Item(913,8).Value = "Brandon Schwab"
$ws.Cells.Item(913,9).Value = 236.5
$ws.Cells.Item(913,10).Value = -1.5
$ws.Cells.Item(913,11).Value = 122
$ws.Cells.Item(913,12).Value = 98
$ws.Cells.Item(913,13).Formula = "=K913+L913"
$ws.Cells.Item(913,14).Formula = "=(L913-K913)*-1"
$ws.Cells.Item(913,15).Value = 1
$ws.Cells.Item(913,16).Formula = "=IF(M913>I913,1,0)"
$ws.Cells.Item(913,17).Formula = "=IF(P913=1,(M913-I913), """")"
$ws.Cells.Item(913,18).Formula = "=IF(M913<I913, 1, 0)"
$ws.Cells.Item(913,19).Formula = "=IF(R913=1,(I913-M913),"""")"
$ws.Cells.Item(913,20).Formula = "=IF(M913=I913,1,0)"

# Row 914: UTA @ DET (Monday 45719)
$ws.Cells.Item(914,1).Value = 45719
$ws.Cells.Item(914,2).Value = "Monday"
$ws.Cells.Item(914,3).Value = "DET"
$ws.Cells.Item(914,4).Value = "UTA"
$ws.Cells.Item(914,5).Value = 0
$ws.Cells.Item(914,6).Value = "Bill Kennedy"
$ws.Cells.Item(914,7).Value = "Tre Maddox"
$ws.Cells.Item(914,8).Value = "Intae Hwang"
$ws.Cells.Item(914,9).Value = 229.5
$ws.Cells.Item(914,10).Value = 9
$ws.Cells.Item(914,11).Value = 134
$ws.Cells.Item(914,12).Value = 106
$ws.Cells.Item(914,13).Formula = "=K914+L914"
$ws.Cells.Item(914,14).Formula = "=(L914-K914)*-1"
$ws.Cells.Item(914,15).Value = 1
$ws.Cells.Item(914,16).Formula = "=IF(M914>I914,1,0)"
$ws.Cells.Item(914,17).Formula = "=IF(P914=1,(M914-I914), """")"
$ws.Cells.Item(914,18).Formula = "=IF(M914<I914, 1, 0)"
$ws.Cells.Item(914,19).Formula = "=IF(R914=1,(I914-M914),"""")"
$ws.Cells.Item(914,20).Formula = "=IF(M914=I914,1,0)"

# Row 915: ORL @ TOR (Tuesday 45720)
$ws.Cells.Item(915,1).Value = 45720
$ws.Cells.Item(915,2).Value = "Tuesday"
$ws.Cells.Item(915,3).Value = "TOR"
$ws.Cells.Item(915,4).Value = "ORL"
$ws.Cells.Item(915,5).Value = 0
$ws.Cells.Item(915,6).Value = "Ben Taylor"
$ws.Cells.Item(915,7).Value = "JT Orr"
$ws.Cells.Item(915,8).Value = "John Conley"
$ws.Cells.Item(915,9).Value = 208.5
$ws.Cells.Item(915,10).Value = -6.5
$ws.Cells.Item(915,11).Value = 114
$ws.Cells.Item(915,12).Value = 113
$ws.Cells.Item(915,13).Formula = "=K915+L915"
$ws.Cells.Item(915,14).Formula = "=(L915-K915)*-1"
$ws.Cells.Item(915,15).Value = 1
$ws.Cells.Item(915,16).Formula = "=IF(M915>I915,1,0)"
$ws.Cells.Item(915,17).Formula = "=IF(P915=1,(M915-I915), """")"
$ws.Cells.Item(915,18).Formula = "=IF(M915<I915, 1, 0)"
$ws.Cells.Item(915,19).Formula = "=IF(R915=1,(I915-M915),"""")"
$ws.Cells.Item(915,20).Formula = "=IF(M915=I915,1,0)"

# Row 916: IND @ HOU (Tuesday 45720)
$ws.Cells.Item(916,1).Value = 45720
$ws.Cells.Item(916,2).Value = "Tuesday"
$ws.Cells.Item(916,3).Value = "HOU"
$ws.Cells.Item(916,4).Value = "IND"
$ws.Cells.Item(916,5).Value = 0
$ws.Cells.Item(916,6).Value = "James Williams"
$ws.Cells.Item(916,7).Value = "Brent Barnaky"
$ws.Cells.Item(916,8).Value = "Brandon Adair"
$ws.Cells.Item(916,9).Value = 229
$ws.Cells.Item(916,10).Value = -4
$ws.Cells.Item(916,11).Value = 102
$ws.Cells.Item(916,12).Value = 115
$ws.Cells.Item(916,13).Formula = "=K916+L916"
$ws.Cells.Item(916,14).Formula = "=(L916-K916)*-1"
$ws.Cells.Item(916,15).Value = 1
$ws.Cells.Item(916,16).Formula = "=IF(M916>I916,1,0)"
$ws.Cells.Item(916,17).Formula = "=IF(P916=1,(M916-I916), """")"
$ws.Cells.Item(916,18).Formula = "=IF(M916<I916, 1, 0)"
$ws.Cells.Item(916,19).Formula = "=IF(R916=1,(I916-M916),"""")"
$ws.Cells.Item(916,20).Formula = "=IF(M916=I916,1,0)"

# Row 917: NYK @ GSW (Tuesday 45720)
$ws.Cells.Item(917,1).Value = 45720
$ws.Cells.Item(917,2).Value = "Tuesday"
$ws.Cells.Item(917,3).Value = "GSW"
$ws.Cells.Item(917,4).Value = "NYK"
$ws.Cells.Item(917,5).Value = 0
$ws.Cells.Item(917,6).Value = "Kevin Scott"
$ws.Cells.Item(917,7).Value = "Natalie Sago"
$ws.Cells.Item(917,8).Value = "CJ Washington"
$ws.Cells.Item(917,9).Value = 232
$ws.Cells.Item(917,10).Value = -3.5
$ws.Cells.Item(917,11).Value = 114
$ws.Cells.Item(917,12).Value = 102
$ws.Cells.Item(917,13).Formula = "=K917+L917"
$ws.Cells.Item(917,14).Formula = "=(L917-K917)*-1"
$ws.Cells.Item(917,15).Value = 1
$ws.Cells.Item(917,16).Formula = "=IF(M917>I917,1,0)"
$ws.Cells.Item(917,17).Formula = "=IF(P917=1,(M917-I917), """")"
$ws.Cells.Item(917,18).Formula = "=IF(M917<I917, 1, 0)"
$ws.Cells.Item(917,19).Formula = "=IF(R917=1,(I917-M917),"""")"
$ws.Cells.Item(917,20).Formula = "=IF(M917=I917,1,0)"

# Row 918: ATL @ MIL (Tuesday 45720)
$ws.Cells.Item(918,1).Value = 45720
$ws.Cells.Item(918,2).Value = "Tuesday"
$ws.Cells.Item(918,3).Value = "MIL"
$ws.Cells.Item(918,4).Value = "ATL"
$ws.Cells.Item(918,5).Value = 0
$ws.Cells.Item(918,6).Value = "Pat Fraher"
$ws.Cells.Item(918,7).Value = "Rodney Mott"
$ws.Cells.Item(918,8).Value = "Robert Hussey"
$ws.Cells.Item(918,9).Value = 244
$ws.Cells.Item(918,10).Value = 4.5
$ws.Cells.Item(918,11).Value = 127
$ws.Cells.Item(918,12).Value = 121
$ws.Cells.Item(918,13).Formula = "=K918+L918"
$ws.Cells.Item(918,14).Formula = "=(L918-K918)*-1"
$ws.Cells.Item(918,15).Value = 1
$ws.Cells.Item(918,16).Formula = "=IF(M918>I918,1,0)"
$ws.Cells.Item(918,17).Formula = "=IF(P918=1,(M918-I918), """")"
$ws.Cells.Item(918,18).Formula = "=IF(M918<I918, 1, 0)"
$ws.Cells.Item(918,19).Formula = "=IF(R918=1,(I918-M918),"""")"
$ws.Cells.Item(918,20).Formula = "=IF(M918=I918,1,0)"

# Row 919: CHI @ CLE (Tuesday 45720)
$ws.Cells.Item(919,1).Value = 45720
$ws.Cells.Item(919,2).Value = "Tuesday"
$ws.Cells.Item(919,3).Value = "CLE"
$ws.Cells.Item(919,4).Value = "CHI"
$ws.Cells.Item(919,5).Value = 0
$ws.Cells.Item(919,6).Value = "Courtney Kirkland"
$ws.Cells.Item(919,7).Value = "Sean Corbin"
$ws.Cells.Item(919,8).Value = "Evan Scott"
$ws.Cells.Item(919,9).Value = 247
$ws.Cells.Item(919,10).Value = 13
$ws.Cells.Item(919,11).Value = 139
$ws.Cells.Item(919,12).Value = 117
$ws.Cells.Item(919,13).Formula = "=K919+L919"
$ws.Cells.Item(919,14).Formula = "=(L919-K919)*-1"
$ws.Cells.Item(919,15).Value = 1
$ws.Cells.Item(919,16).Formula = "=IF(M919>I919,1,0)"
$ws.Cells.Item(919,17).Formula = "=IF(P919=1,(M919-I919), """")"
$ws.Cells.Item(919,18).Formula = "=IF(M919<I919, 1, 0)"
$ws.Cells.Item(919,19).Formula = "=IF(R919=1,(I919-M919),"""")"
$ws.Cells.Item(919,20).Formula = "=IF(M919=I919,1,0)"

# Row 920: MIN @ PHI (Tuesday 45720)
$ws.Cells.Item(920,1).Value = 45720
$ws.Cells.Item(920,2).Value = "Tuesday"
$ws.Cells.Item(920,3).Value = "PHI"
$ws.Cells.Item(920,4).Value = "MIN"
$ws.Cells.Item(920,5).Value = 0
$ws.Cells.Item(920,6).Value = "Josh Tiven"
$ws.Cells.Item(920,7).Value = "Mousa Dagher"
$ws.Cells.Item(920,8).Value = "ShaRae Mitchell"
$ws.Cells.Item(920,9).Value = 221.5
$ws.Cells.Item(920,10).Value = -11.5
$ws.Cells.Item(920,11).Value = 112
$ws.Cells.Item(920,12).Value = 126
$ws.Cells.Item(920,13).Formula = "=K920+L920"
$ws.Cells.Item(920,14).Formula = "=(L920-K920)*-1"
$ws.Cells.Item(920,15).Value = 1
$ws.Cells.Item(920,16).Formula = "=IF(M920>I920,1,0)"
$ws.Cells.Item(920,17).Formula = "=IF(P920=1,(M920-I920), """")"
$ws.Cells.Item(920,18).Formula = "=IF(M920<I920, 1, 0)"
$ws.Cells.Item(920,19).Formula = "=IF(R920=1,(I920-M920),"""")"
$ws.Cells.Item(920,20).Formula = "=IF(M920=I920,1,0)"

# Row 921: SAS @ BKN (Tuesday 45720)
$ws.Cells.Item(921,1).Value = 45720
$ws.Cells.Item(921,2).Value = "Tuesday"
$ws.Cells.Item(921,3).Value = "BKN"
$ws.Cells.Item(921,4).Value = "SAS"
$ws.Cells.Item(921,5).Value = 0
$ws.Cells.Item(921,6).Value = "David Guthrie"
$ws.Cells.Item(921,7).Value = "Nick Buchert"
$ws.Cells.Item(921,8).Value = "Simone Jelks"
$ws.Cells.Item(921,9).Value = 221.5
$ws.Cells.Item(921,10).Value = -5.5
$ws.Cells.Item(921,11).Value = 113
$ws.Cells.Item(921,12).Value = 127
$ws.Cells.Item(921,13).Formula = "=K921+L921"
$ws.Cells.Item(921,14).Formula = "=(L921-K921)*-1"
$ws.Cells.Item(921,15).Value = 1
$ws.Cells.Item(921,16).Formula = "=IF(M921>I921,1,0)"
$ws.Cells.Item(921,17).Formula = "=IF(P921=1,(M921-I921), """")"
$ws.Cells.Item(921,18).Formula = "=IF(M921<I921, 1, 0)"
$ws.Cells.Item(921,19).Formula = "=IF(R921=1,(I921-M921),"""")"
$ws.Cells.Item(921,20).Formula = "=IF(M921=I921,1,0)"

# Row 922: PHX @ LAC (Tuesday 45720)
$ws.Cells.Item(922,1).Value = 45720
$ws.Cells.Item(922,2).Value = "Tuesday"
$ws.Cells.Item(922,3).Value = "LAC"
$ws.Cells.Item(922,4).Value = "PHX"
$ws.Cells.Item(922,5).Value = 0
$ws.Cells.Item(922,6).Value = "Brian Forte"
$ws.Cells.Item(922,7).Value = "JB DeRosa"
$ws.Cells.Item(922,8).Value = "Nate Green"
$ws.Cells.Item(922,9).Value = 228.5
$ws.Cells.Item(922,10).Value = -1
$ws.Cells.Item(922,11).Value = 117
$ws.Cells.Item(922,12).Value = 119
$ws.Cells.Item(922,13).Formula = "=K922+L922"
$ws.Cells.Item(922,14).Formula = "=(L922-K922)*-1"
$ws.Cells.Item(922,15).Value = 1
$ws.Cells.Item(922,16).Formula = "=IF(M922>I922,1,0)"
$ws.Cells.Item(922,17).Formula = "=IF(P922=1,(M922-I922), """")"
$ws.Cells.Item(922,18).Formula = "=IF(M922<I922, 1, 0)"
$ws.Cells.Item(922,19).Formula = "=IF(R922=1,(I922-M922),"""")"
$ws.Cells.Item(922,20).Formula = "=IF(M922=I922,1,0)"

# Row 923: LAL @ NOP (Tuesday 45720)
$ws.Cells.Item(923,1).Value = 45720
$ws.Cells.Item(923,2).Value = "Tuesday"
$ws.Cells.Item(923,3).Value = "NOP"
$ws.Cells.Item(923,4).Value = "LAL"
$ws.Cells.Item(923,5).Value = 0
$ws.Cells.Item(923,6).Value = "Karl Lane"
$ws.Cells.Item(923,7).Value = "Kevin Cutler"
$ws.Cells.Item(923,8).Value = "Matt Kallio"
$ws.Cells.Item(923,9).Value = 234.5
$ws.Cells.Item(923,10).Value = -8
$ws.Cells.Item(923,11).Value = 115
$ws.Cells.Item(923,12).Value = 136
$ws.Cells.Item(923,13).Formula = "=K923+L923"
$ws.Cells.Item(923,14).Formula = "=(L923-K923)*-1"
$ws.Cells.Item(923,15).Value = 1
$ws.Cells.Item(923,16).Formula = "=IF(M923>I923,1,0)"
$ws.Cells.Item(923,17).Formula = "=IF(P923=1,(M923-I923), """")"
$ws.Cells.Item(923,18).Formula = "=IF(M923<I923, 1, 0)"
$ws.Cells.Item(923,19).Formula = "=IF(R923=1,(I923-M923),"""")"
$ws.Cells.Item(923,20).Formula = "=IF(M923=I923,1,0)"

# --- Step 4: view state - select next empty row and scroll the frozen pane down ---
[void]$ws.Range("A924").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 883
